$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Insert a new column at N (14th column) on the hidden "#system" sheet.
#    This shifts the existing category columns N..AC (macro..xml) one to the
#    right, becoming O..AD, and makes room for the new "localdb" category.
# ---------------------------------------------------------------------------
$ws.Columns(14).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new column N with the "localdb" category: a header plus
#    the six new command names (alphabetically ordered), mirroring the
#    layout used by every other category column on this sheet.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 14).Value = "localdb"
$ws.Cells.Item(2, 14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value = "purge(var)"
$ws.Cells.Item(7, 14).Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------------
# 3. The "target" column (A) holds every category name in alphabetical
#    order, independent of the other columns. "localdb" sorts right before
#    "macro", so insert a single cell there (shifting macro..xml down one
#    row within column A only) and write the new category name in.
# ---------------------------------------------------------------------------
$ws.Range("A14").Insert()
$ws.Cells.Item(14, 1).Value = "localdb"

# ---------------------------------------------------------------------------
# 4. Defined names do not auto-adjust to the column insert, so update every
#    name whose range shifted right by one column, and add the new
#    "localdb" defined name.
# ---------------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")

Write-Host "localdb command category added"
